$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers: extend Accuracy / F-measure pattern across D1:G1
$ws.Range("D1").Value = "Accuracy"
$ws.Range("E1").Value = "F-measure"
$ws.Range("F1").Value = "Accuracy"
$ws.Range("G1").Value = "F-measure"

# Row 2 (Mean)
$ws.Range("D2").Value = 82.6
$ws.Range("E2").Value = 0.832
$ws.Range("F2").Value = 89.13
$ws.Range("G2").Value = 0.889

# Row 3 (Min)
$ws.Range("D3").Value = 69.56
$ws.Range("E3").Value = 0.72
$ws.Range("F3").Value = 80.43
$ws.Range("G3").Value = 0.808

# Row 4 (Max)
$ws.Range("D4").Value = 82.6
$ws.Range("E4").Value = 0.789
$ws.Range("F4").Value = 84.78
$ws.Range("G4").Value = 0.844

# Row 5 (Mode)
$ws.Range("D5").Value = 78.26
$ws.Range("E5").Value = 0.772
$ws.Range("F5").Value = 84.78
$ws.Range("G5").Value = 0.836

# Row 6 (Median)
$ws.Range("D6").Value = 80.43
$ws.Range("E6").Value = 0.717
$ws.Range("F6").Value = 82.6
$ws.Range("G6").Value = 0.826

# Row 7-8 labels
$ws.Range("D7").Value = "Decision Tree"
$ws.Range("F7").Value = "Random Forest"
$ws.Range("D8").Value = "70-30-split"
$ws.Range("F8").Value = "70-30-split"

# View state: zoom out and move the selection, matching the author's session
$excel.ActiveWindow.Zoom = 200
[void]$ws.Range("F12").Select()
